$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new values
$ws.Range("A2").Value = "76513680-6/0"
$ws.Range("B2").Value = "CFINHRFLA"
$ws.Range("C2").Value = "L"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 16960.8644
$ws.Range("F2").Value = "21/03/2025"
$ws.Range("G2").Value = "21/03/2025"
$ws.Range("H2").Value = 33922
$ws.Range("I2").Value = 33922
$ws.Range("J2").Value = 0

# Remove rows 3 and 4 entirely
$ws.Rows("3:4").Delete()
